# Weekly refresh of fruit/vegetable (Haba) price rows.
# Each data row (2-8, 10-16) gets new Fecha/Volumen/Precio values - a
# re-shuffle of the weekly price observations already present elsewhere
# in the sheet. Row 9 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(Fecha, Volumen, PrecioMinimo, PrecioMaximo, PrecioPromedio, Origen, PrecioKg)
$updates = @{
    2  = @(44372, 500, 13000, 14000, 13500, "Provincia del Elquí", 540)
    3  = @(44690, 400, 17000, 18000, 17500, "Provincia del Elquí", 700)
    4  = @(44370, 520, 13000, 14000, 13500, "Provincia del Elquí", 540)
    5  = @(44376, 400, 12000, 13000, 12500, "Provincia del Elquí", 500)
    6  = @(44484, 400,  9000, 10000,  9500, "Provincia del Elquí", 380)
    7  = @(44473, 500,  8500,  9000,  8750, "Provincia del Elquí", 350)
    8  = @(44386, 500, 11000, 12000, 11500, "Provincia del Elquí", 460)
    10 = @(44384, 560, 11500, 12000, 11750, "Provincia del Elquí", 470)
    11 = @(44694, 480, 17500, 18000, 17750, "Provincia del Elquí", 710)
    12 = @(44466, 400,  9500, 10000,  9750, "Provincia del Elquí", 390)
    13 = @(44356, 500, 13000, 14000, 13500, "Provincia de Limarí", 540)
    14 = @(44316, 300, 16000, 17000, 16500, "Provincia del Elquí", 660)
    15 = @(44425, 400, 11500, 12000, 11750, "Provincia del Elquí", 470)
    16 = @(44377, 520, 12500, 13000, 12750, "Provincia del Elquí", 510)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K - Precio mínimo
    $ws.Cells.Item($row, 12).Value = $vals[3]   # L - Precio máximo
    $ws.Cells.Item($row, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals[5]   # O - Origen
    $ws.Cells.Item($row, 16).Value = $vals[6]   # P - Precio $/Kg
}
